$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("H2").Value = 578
$ws.Range("J2").Value = 6741
$ws.Range("L2").Value = 1879
$ws.Range("M2").Value = 128
$ws.Range("N2").Value = 1178
$ws.Range("O2").Value = 6
$ws.Range("P2").Value = 38
$ws.Range("Q2").Value = 17
$ws.Range("R2").Value = 97
$ws.Range("S2").Value = 751
$ws.Range("T2").Value = 1148
$ws.Range("U2").Value = 114
$ws.Range("V2").Value = 10618
$ws.Range("W2").Value = 4
$ws.Range("X2").Value = 10264
$ws.Range("Y2").Value = 11
$ws.Range("Z2").Value = 154
$ws.Range("AA2").Value = 57
